# update 03 juli 2023
# Change the wording of the question in cell D5 ("Apa pendapat anda tentang
# 2 Petrus 1:5-7 ? Jelaskan " -> "... Mohon dijelaskan ") and leave the
# cursor/selection on D14, matching the author's last saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Apa pendapat anda tentang 2 Petrus 1:5-7 ? Mohon dijelaskan "

$ws.Range("D14").Select()
